$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4128
$ws.Range("E2").Value = 132
$ws.Range("F2").Value = 181
$ws.Range("G2").Value = -30
$ws.Range("H2").Value = -24
$ws.Range("I2").Value = -24
$ws.Range("K2").Value = 5892
$ws.Range("L2").Value = 3692
$ws.Range("M2").Value = 2200
$ws.Range("N2").Value = 2200
$ws.Range("P2").Value = 415
$ws.Range("Q2").Value = 73
$ws.Range("R2").Value = 108
$ws.Range("S2").Value = -141
$ws.Range("T2").Value = 69
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 2775
$ws.Range("W2").Value = 3.21
$ws.Range("X2").Value = -0.57
$ws.Range("Y2").Value = -1.12
$ws.Range("Z2").Value = -0.4
$ws.Range("AA2").Value = 167.81
$ws.Range("AB2").Value = 407.31
$ws.Range("AC2").Value = -124
$ws.Range("AD2").Value = -101.59
$ws.Range("AE2").Value = 11121
$ws.Range("AF2").Value = 1.13
$ws.Range("AG2").Value = 104
$ws.Range("AH2").Value = 0.82
$ws.Range("AI2").Value = -115.04
$ws.Range("AJ2").Value = 18139164
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 4344
$ws.Range("E3").Value = 215
$ws.Range("F3").Value = 217
$ws.Range("G3").Value = 21
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 20
$ws.Range("K3").Value = 6135
$ws.Range("L3").Value = 3881
$ws.Range("M3").Value = 2254
$ws.Range("N3").Value = 2254
$ws.Range("P3").Value = 452
$ws.Range("Q3").Value = -42
$ws.Range("R3").Value = -73
$ws.Range("S3").Value = 147
$ws.Range("T3").Value = 56
$ws.Range("U3").Value = -97
$ws.Range("V3").Value = 2905
$ws.Range("W3").Value = 4.96
$ws.Range("X3").Value = 0.45
$ws.Range("Y3").Value = 0.88
$ws.Range("Z3").Value = 0.33
$ws.Range("AA3").Value = 172.17
$ws.Range("AB3").Value = 377.71
$ws.Range("AC3").Value = 95
$ws.Range("AD3").Value = 352.55
$ws.Range("AE3").Value = 10889
$ws.Range("AF3").Value = 3.08
$ws.Range("AG3").Value = 131
$ws.Range("AH3").Value = 0.39
$ws.Range("AI3").Value = 138.33
$ws.Range("AJ3").Value = 20473468
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 4675
$ws.Range("E4").Value = 245
$ws.Range("F4").Value = 245
$ws.Range("G4").Value = 71
$ws.Range("H4").Value = -109
$ws.Range("I4").Value = -109
$ws.Range("K4").Value = 6511
$ws.Range("L4").Value = 3748
$ws.Range("M4").Value = 2763
$ws.Range("N4").Value = 2763
$ws.Range("P4").Value = 501
$ws.Range("Q4").Value = 67
$ws.Range("R4").Value = -88
$ws.Range("S4").Value = 338
$ws.Range("T4").Value = 81
$ws.Range("U4").Value = -14
$ws.Range("V4").Value = 2512
$ws.Range("W4").Value = 5.25
$ws.Range("X4").Value = -2.34
$ws.Range("Y4").Value = -4.36
$ws.Range("Z4").Value = -1.73
$ws.Range("AA4").Value = 135.63
$ws.Range("AB4").Value = 432.17
$ws.Range("AC4").Value = -514
$ws.Range("AD4").Value = -96.02
$ws.Range("AE4").Value = 12270
$ws.Range("AF4").Value = 4.02
$ws.Range("AG4").Value = 229
$ws.Range("AH4").Value = 0.46
$ws.Range("AI4").Value = -47.51
$ws.Range("AJ4").Value = 20473468
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 5029
$ws.Range("E5").Value = 217
$ws.Range("F5").Value = 217
$ws.Range("G5").Value = -13
$ws.Range("H5").Value = -7
$ws.Range("I5").Value = -7
$ws.Range("K5").Value = 6309
$ws.Range("L5").Value = 3601
$ws.Range("M5").Value = 2708
$ws.Range("N5").Value = 2708
$ws.Range("P5").Value = 512
$ws.Range("Q5").Value = 334
$ws.Range("R5").Value = -99
$ws.Range("S5").Value = -279
$ws.Range("T5").Value = 73
$ws.Range("U5").Value = 262
$ws.Range("V5").Value = 2340
$ws.Range("W5").Value = 4.32
$ws.Range("X5").Value = -0.15
$ws.Range("Y5").Value = -0.27
$ws.Range("Z5").Value = -0.12
$ws.Range("AA5").Value = 132.99
$ws.Range("AB5").Value = 410.94
$ws.Range("AC5").Value = -33
$ws.Range("AD5").Value = -1312.46
$ws.Range("AE5").Value = 12016
$ws.Range("AF5").Value = 3.6
$ws.Range("AG5").Value = 257
$ws.Range("AH5").Value = 0.59
$ws.Range("AI5").Value = -784.01
$ws.Range("AJ5").Value = 20510014
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 5372
$ws.Range("E6").Value = 216
$ws.Range("F6").Value = 216
$ws.Range("G6").Value = 102
$ws.Range("H6").Value = 81
$ws.Range("I6").Value = 81
$ws.Range("K6").Value = 5938
$ws.Range("L6").Value = 3337
$ws.Range("M6").Value = 2601
$ws.Range("N6").Value = 2601
$ws.Range("P6").Value = 528
$ws.Range("Q6").Value = 521
$ws.Range("R6").Value = -139
$ws.Range("S6").Value = -412
$ws.Range("T6").Value = 77
$ws.Range("U6").Value = 444
$ws.Range("V6").Value = 2028
$ws.Range("W6").Value = 4.02
$ws.Range("X6").Value = 1.5
$ws.Range("Y6").Value = 3.03
$ws.Range("Z6").Value = 1.32
$ws.Range("AA6").Value = 128.28
$ws.Range("AB6").Value = 411.41
$ws.Range("AC6").Value = 356
$ws.Range("AD6").Value = 111.24
$ws.Range("AE6").Value = 11616
$ws.Range("AF6").Value = 3.41
$ws.Range("AG6").Value = 310
$ws.Range("AH6").Value = 0.78
$ws.Range("AI6").Value = 86.31999999999999
$ws.Range("AJ6").Value = 20795121

# Row 7: clear all numeric data, keep A/B/C
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all numeric data, keep A/B/C
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all numeric data, keep A/B/C
$ws.Range("D9:AJ9").ClearContents()

